# Daily attendance processing - 2025-12-25 19:28:34
# Normalizes the "Recorded By" (column G) cell values so that the
# "System"/"system" token appears first in the comma-separated list,
# matching the canonical ordering used by the attendance export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of exact "before" -> "after" strings for the Recorded By column.
$map = @{
    "System, backup@backdoor.com, system" = "System, system, backup@backdoor.com"
    "admin@admin.com, System"              = "System, admin@admin.com"
    "dnasr281@gmail.com, System"           = "System, dnasr281@gmail.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value()
    if ($map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
